$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Elimina los periodos de mora (EC) anteriores y se agregan nuevos (base de datos),
# invirtiendo el orden de los periodos existentes.
$ws.Range("E16").Value = "2306"
$ws.Range("E17").Value = "2305"
$ws.Range("E18").Value = "2304"
$ws.Range("E19").Value = "2303"
